# Updates cryptos list values (price/volume) to reflect the latest scrape,
# including a row swap for WrappedBTC/ShibaInu (rows 15/16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''70.130.50'
$ws.Range("E2").Value = '  -3.08%  '
$ws.Range("D3").Value = '''2.517.75'
$ws.Range("E3").Value = '  -5.48%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''574.09'
$ws.Range("E5").Value = '  -4.02%  '
$ws.Range("D6").Value = '''169.53'
$ws.Range("E6").Value = '  -3.31%  '
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("D8").Value = '''0.509'
$ws.Range("E8").Value = '  -2.87%  '
$ws.Range("D9").Value = '''2.515.61'
$ws.Range("E9").Value = '  -5.52%  '
$ws.Range("E10").Value = '  -5.18%  '
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("E12").Value = '  -3.71%  '
$ws.Range("D13").Value = '''4.80'
$ws.Range("E13").Value = '  -4.13%  '
$ws.Range("D14").Value = '''2.990.15'
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '''70.029.59'
$ws.Range("E15").Value = '  -2.92%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '''0.0000179'
$ws.Range("E16").Value = '  -2.96%  '
$ws.Range("D17").Value = '''24.89'
$ws.Range("E17").Value = '  -5.25%  '
$ws.Range("D18").Value = '''2.525.47'
$ws.Range("E18").Value = '  -5.26%  '
$ws.Range("D19").Value = '''11.40'
$ws.Range("E19").Value = '  -6.85%  '
$ws.Range("E20").Value = '  -8.48%  '
$ws.Range("D21").Value = '''353.41'
$ws.Range("E21").Value = '  -4.66%  '
$ws.Range("E22").Value = '  -5.86%  '
$ws.Range("E23").Value = '  -3.25%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = '''68.87'
$ws.Range("E25").Value = '  -4.33%  '
$ws.Range("E26").Value = '  -5.80%  '
$ws.Range("D27").Value = '''9.22'
$ws.Range("E27").Value = '  -5.53%  '
$ws.Range("D28").Value = '''2.651.52'
$ws.Range("E28").Value = '  -5.26%  '
$ws.Range("D29").Value = '''0.998'
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("D30").Value = '''0.0₃0909'
$ws.Range("E30").Value = '  -6.07%  '
$ws.Range("D31").Value = '''7.82'
$ws.Range("E31").Value = '  -3.15%  '
$ws.Range("D32").Value = '''479.37'
$ws.Range("E32").Value = '  -4.24%  '
$ws.Range("E33").Value = '  -0.66%  '
$ws.Range("E34").Value = '  -4.05%  '
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").Value = '''156.96'
$ws.Range("E36").Value = '  -3.68%  '
$ws.Range("D37").Value = '''0.116'
$ws.Range("E37").Value = '  +3.09%  '
$ws.Range("D38").Value = '''18.84'
$ws.Range("E38").Value = '  -0.67%  '
$ws.Range("D39").Value = '''18.55'
$ws.Range("E39").Value = '  -5.03%  '
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("E41").Value = '  -5.92%  '
$ws.Range("D42").Value = '''1.64'
$ws.Range("E42").Value = '  -7.40%  '
$ws.Range("E43").Value = '  -4.13%  '
$ws.Range("E44").Value = '  -5.56%  '
$ws.Range("D45").Value = '''2.37'
$ws.Range("E45").Value = '  -6.96%  '
$ws.Range("D46").Value = '''38.28'
$ws.Range("E46").Value = '  -3.15%  '
$ws.Range("D47").Value = '''141.93'
$ws.Range("E47").Value = '  -9.31%  '
$ws.Range("E48").Value = '  -5.86%  '
$ws.Range("E49").Value = '  -6.37%  '
$ws.Range("E50").Value = '  -7.14%  '
$ws.Range("D51").Value = '''0.596'
$ws.Range("E51").Value = '  -1.62%  '
